$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Direct text / link updates (not ambiguous with numeric parsing) ---
$ws.Range("D2").Value = '45.149.47'
$ws.Range("E2").Value = '  -3.44%  '
$ws.Range("D3").Value = '2.380.87'
$ws.Range("E3").Value = '  +4.82%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("E5").Value = '  -3.43%  '
$ws.Range("E6").Value = '  -6.43%  '
$ws.Range("E7").Value = '  -1.62%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -3.62%  '
$ws.Range("E10").Value = '  -5.42%  '
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("E12").Value = '  -4.15%  '
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("D14").Value = '2.745.20'
$ws.Range("E14").Value = '  +4.89%  '
$ws.Range("D15").Value = '2.380.31'
$ws.Range("E15").Value = '  +5.00%  '
$ws.Range("E16").Value = '  +1.25%  '
$ws.Range("E17").Value = '  +2.10%  '
$ws.Range("D18").Value = '45.111.03'
$ws.Range("E18").Value = '  -3.39%  '
$ws.Range("E19").Value = '  -6.98%  '
$ws.Range("E20").Value = '  -0.77%  '
$ws.Range("E21").Value = '  +2.23%  '
$ws.Range("E22").Value = '  +0.97%  '
$ws.Range("E23").Value = '  -4.66%  '
$ws.Range("E24").Value = '  -3.95%  '
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("E26").Value = '  -0.86%  '
$ws.Range("E27").Value = '  -1.57%  '
$ws.Range("E28").Value = '  -13.15%  '
$ws.Range("E29").Value = '  -2.81%  '
$ws.Range("E30").Value = '  +13.85%  '
$ws.Range("E31").Value = '  +4.86%  '
$ws.Range("E32").Value = '  -2.99%  '
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("E34").Value = '  -2.95%  '
$ws.Range("E35").Value = '  -3.18%  '
$ws.Range("E36").Value = '  -3.61%  '
$ws.Range("E37").Value = '  +11.24%  '
$ws.Range("E38").Value = '  -2.41%  '
$ws.Range("E39").Value = '  -9.12%  '
$ws.Range("E40").Value = '  -5.70%  '
$ws.Range("E41").Value = '  -2.06%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.965.27'
$ws.Range("E42").Value = '  +8.55%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("E43").Value = '  -4.06%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("E45").Value = '  -2.18%  '
$ws.Range("E46").Value = '  -13.87%  '
$ws.Range("E47").Value = '  +7.35%  '
$ws.Range("E48").Value = '  +14.72%  '
$ws.Range("E49").Value = '  +4.49%  '
$ws.Range("D50").Value = '2.615.65'
$ws.Range("E50").Value = '  +4.87%  '
$ws.Range("E51").Value = '  -5.03%  '

# --- Numeric-looking price updates: these must stay text cells, matching the
# workbook's original inline-string storage, so we quote-prefix them to force
# text interpretation, then restore the default "Normal" style so no stray
# number-format / quote-prefix styling is left behind on the cell.
$ws.Range("D5").Value = "'292.66"
$ws.Range("D5").Style = 'Normal'
$ws.Range("D6").Value = "'93.79"
$ws.Range("D6").Style = 'Normal'
$ws.Range("D7").Value = "'0.553"
$ws.Range("D7").Style = 'Normal'
$ws.Range("D9").Value = "'0.495"
$ws.Range("D9").Style = 'Normal'
$ws.Range("D10").Value = "'33.82"
$ws.Range("D10").Style = 'Normal'
$ws.Range("D11").Value = "'0.0772"
$ws.Range("D11").Style = 'Normal'
$ws.Range("D12").Value = "'6.91"
$ws.Range("D12").Style = 'Normal'
$ws.Range("D16").Value = "'13.86"
$ws.Range("D16").Style = 'Normal'
$ws.Range("D17").Value = "'0.817"
$ws.Range("D17").Style = 'Normal'
$ws.Range("D19").Value = "'12.26"
$ws.Range("D19").Style = 'Normal'
$ws.Range("D21").Value = "'6.05"
$ws.Range("D21").Style = 'Normal'
$ws.Range("D22").Value = "'65.96"
$ws.Range("D22").Style = 'Normal'
$ws.Range("D23").Value = "'237.34"
$ws.Range("D23").Style = 'Normal'
$ws.Range("D27").Value = "'2.21"
$ws.Range("D27").Style = 'Normal'
$ws.Range("D28").Value = "'37.08"
$ws.Range("D28").Style = 'Normal'
$ws.Range("D29").Value = "'9.46"
$ws.Range("D29").Style = 'Normal'
$ws.Range("D30").Value = "'3.77"
$ws.Range("D30").Style = 'Normal'
$ws.Range("D31").Value = "'20.87"
$ws.Range("D31").Style = 'Normal'
$ws.Range("D33").Value = "'146.76"
$ws.Range("D33").Style = 'Normal'
$ws.Range("D34").Value = "'5.35"
$ws.Range("D34").Style = 'Normal'
$ws.Range("D35").Value = "'0.0753"
$ws.Range("D35").Style = 'Normal'
$ws.Range("D36").Value = "'0.111"
$ws.Range("D36").Style = 'Normal'
$ws.Range("D37").Value = "'1.93"
$ws.Range("D37").Style = 'Normal'
$ws.Range("D38").Value = "'0.113"
$ws.Range("D38").Style = 'Normal'
$ws.Range("D39").Value = "'14.61"
$ws.Range("D39").Style = 'Normal'
$ws.Range("D40").Value = "'3.69"
$ws.Range("D40").Style = 'Normal'
$ws.Range("D43").Value = "'3.13"
$ws.Range("D43").Style = 'Normal'
$ws.Range("D45").Value = "'88.73"
$ws.Range("D45").Style = 'Normal'
$ws.Range("D46").Value = "'1.71"
$ws.Range("D46").Style = 'Normal'
$ws.Range("D47").Value = "'8.40"
$ws.Range("D47").Style = 'Normal'
$ws.Range("D48").Value = "'14.87"
$ws.Range("D48").Style = 'Normal'
$ws.Range("D49").Value = "'98.80"
$ws.Range("D49").Style = 'Normal'
$ws.Range("D51").Value = "'0.180"
$ws.Range("D51").Style = 'Normal'
